# Update "想去人数" (interested-attendee count) figures that changed between
# the previous data pull and the refreshed gh-pages data pull (commit 456a3b4).
#
# Sheet "展览" (Exhibitions) and the aggregated "全部类型" (All types) sheet
# both list the same events, so each numeric change needs to be applied in
# both places. Sheet "演出" (Performances) only needs a single cell updated.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 5527   # was 5519
$wsExhibit.Range("F5").Value = 309    # was 308
$wsExhibit.Range("F7").Value = 37     # was 33

# 演出 (Performances)
$wsShows = $wb.Worksheets.Item("演出")
$wsShows.Range("F3").Value = 18       # was 17

# 全部类型 (All types) - aggregated view of every event
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 5527       # was 5519
$wsAll.Range("F5").Value = 309        # was 308
$wsAll.Range("F7").Value = 37         # was 33
$wsAll.Range("F13").Value = 18        # was 17

$wb.Save()
